# Auto-generated Excel COM-interop script
# Applies updated market-price data values per the commit diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 742686.0600000001
$ws.Range("J17").Value = 742686.0600000001
$ws.Range("L17").Value = 2228058.18
$ws.Range("N17").Value = -2228394.18
$ws.Range("H28").Value = 931.6667
$ws.Range("I28").Value = 359.84616
$ws.Range("K28").Value = 359.84616
$ws.Range("M28").Value = 125.15384
$ws.Range("H32").Value = 4273.5
$ws.Range("J32").Value = 4273.5
$ws.Range("L32").Value = 4273.5
$ws.Range("N32").Value = -4925.5
$ws.Range("H33").Value = 402.2
$ws.Range("I33").Value = 252.75
$ws.Range("K33").Value = 252.75
$ws.Range("M33").Value = -23.75
$ws.Range("H40").Value = 2195.5
$ws.Range("I40").Value = 2234.2
$ws.Range("K40").Value = 2234.2
$ws.Range("M40").Value = -2059.2
$ws.Range("H48").Value = 5019
$ws.Range("J48").Value = 5019
$ws.Range("L48").Value = 15057
$ws.Range("N48").Value = -15641
$ws.Range("H56").Value = 5019
$ws.Range("J56").Value = 5019
$ws.Range("L56").Value = 15057
$ws.Range("N56").Value = -16125
$ws.Range("H62").Value = 6061.2
$ws.Range("I62").Value = 2700
$ws.Range("J62").Value = 6434.6665
$ws.Range("K62").Value = 2700
$ws.Range("L62").Value = 6434.6665
$ws.Range("M62").Value = -2076
$ws.Range("N62").Value = -7682.6665
$ws.Range("H64").Value = 12333.167
$ws.Range("I64").Value = 12399.8
$ws.Range("K64").Value = 12399.8
$ws.Range("M64").Value = -12151.8
$ws.Range("H65").Value = 6061.2
$ws.Range("I65").Value = 2700
$ws.Range("J65").Value = 6434.6665
$ws.Range("K65").Value = 13500
$ws.Range("L65").Value = 32173.3325
$ws.Range("M65").Value = -10380
$ws.Range("N65").Value = -38413.3325
$ws.Range("H67").Value = 12333.167
$ws.Range("I67").Value = 12399.8
$ws.Range("K67").Value = 12399.8
$ws.Range("M67").Value = -11541.8
$ws.Range("H70").Value = 1737.25
$ws.Range("J70").Value = 1737.25
$ws.Range("L70").Value = 5211.75
$ws.Range("N70").Value = -5751.75
$ws.Range("H73").Value = 1737.25
$ws.Range("J73").Value = 1737.25
$ws.Range("L73").Value = 5211.75
$ws.Range("N73").Value = -7083.75
$ws.Range("H76").Value = 9979.049999999999
$ws.Range("I76").Value = 14732.333
$ws.Range("K76").Value = 14732.333
$ws.Range("M76").Value = -14417.333
$ws.Range("H79").Value = 9979.049999999999
$ws.Range("I79").Value = 14732.333
$ws.Range("K79").Value = 14732.333
$ws.Range("M79").Value = -13640.333
$ws.Range("H99").Value = 863.2941
$ws.Range("I99").Value = 523.2857
$ws.Range("J99").Value = 2450
$ws.Range("K99").Value = 1569.8571
$ws.Range("L99").Value = 7350
$ws.Range("M99").Value = -71.85710000000017
$ws.Range("N99").Value = -10346
$ws.Range("H106").Value = 1272.5
$ws.Range("I106").Value = 1272.5
$ws.Range("K106").Value = 1272.5
$ws.Range("M106").Value = -641.5
$ws.Range("H111").Value = 11907084
$ws.Range("I111").Value = 13891272
$ws.Range("K111").Value = 41673816
$ws.Range("M111").Value = -41670749
$ws.Range("H129").Value = 27511.875
$ws.Range("I129").Value = 1150
$ws.Range("K129").Value = 3450
$ws.Range("M129").Value = 1550
$ws.Range("H134").Value = 71616.836
$ws.Range("I134").Value = 39700
$ws.Range("J134").Value = 78000.2
$ws.Range("K134").Value = 39700
$ws.Range("L134").Value = 78000.2
$ws.Range("M134").Value = -34630
$ws.Range("N134").Value = -88140.2
$ws.Range("H136").Value = 64000
$ws.Range("J136").Value = 64000
$ws.Range("L136").Value = 64000
$ws.Range("N136").Value = -74200
$ws.Range("H137").Value = 1867.6744
$ws.Range("J137").Value = 1718.4117
$ws.Range("L137").Value = 5155.2351
$ws.Range("N137").Value = -10255.2351
$ws.Range("H138").Value = 47635104
$ws.Range("I138").Value = 2332
$ws.Range("K138").Value = 6996
$ws.Range("M138").Value = -1856
$ws.Range("H139").Value = 134410.78
$ws.Range("J139").Value = 146249.62
$ws.Range("L139").Value = 146249.62
$ws.Range("N139").Value = -156529.62
$ws.Range("H140").Value = 64849.5
$ws.Range("I140").Value = 39700
$ws.Range("J140").Value = 89999
$ws.Range("K140").Value = 39700
$ws.Range("L140").Value = 89999
$ws.Range("M140").Value = -34520
$ws.Range("N140").Value = -100359

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H23").Value = 16331.667
$ws.Range("I23").Value = 16331.667
$ws.Range("K23").Value = 16331.667
$ws.Range("M23").Value = -16072.667
$ws.Range("H37").Value = 20000
$ws.Range("I37").Value = 20000
$ws.Range("K37").Value = 20000
$ws.Range("M37").Value = -19727
$ws.Range("H61").Value = 11909338
$ws.Range("I61").Value = 13162164
$ws.Range("K61").Value = 13162164
$ws.Range("M61").Value = -13161952
$ws.Range("H74").Value = 24418238
$ws.Range("I74").Value = 25670404
$ws.Range("K74").Value = 25670404
$ws.Range("M74").Value = -25669530
$ws.Range("H77").Value = 24418238
$ws.Range("I77").Value = 25670404
$ws.Range("K77").Value = 128352020
$ws.Range("M77").Value = -128347652
$ws.Range("H122").Value = 2897.1052
$ws.Range("I122").Value = 2193.0952
$ws.Range("J122").Value = 3766.7646
$ws.Range("K122").Value = 6579.285600000001
$ws.Range("L122").Value = 11300.2938
$ws.Range("M122").Value = -4129.285600000001
$ws.Range("N122").Value = -16200.2938
$ws.Range("H136").Value = 11909338
$ws.Range("I136").Value = 13162164
$ws.Range("K136").Value = 39486492
$ws.Range("M136").Value = -39483942

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2343.75
$ws.Range("I20").Value = 2900
$ws.Range("J20").Value = 2158.3333
$ws.Range("K20").Value = 2900
$ws.Range("L20").Value = 2158.3333
$ws.Range("M20").Value = -2653
$ws.Range("N20").Value = -2652.3333
$ws.Range("H94").Value = 803.0714
$ws.Range("I94").Value = 269.6
$ws.Range("K94").Value = 269.6
$ws.Range("M94").Value = 181.4
$ws.Range("H99").Value = 5477.7144
$ws.Range("I99").Value = 4083.25
$ws.Range("K99").Value = 4083.25
$ws.Range("M99").Value = -2585.25
$ws.Range("H107").Value = 1898.4546
$ws.Range("I107").Value = 1724.96
$ws.Range("K107").Value = 1724.96
$ws.Range("M107").Value = 195.04

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 350.55554
$ws.Range("I7").Value = 231.7
$ws.Range("J7").Value = 499.125
$ws.Range("K7").Value = 231.7
$ws.Range("L7").Value = 499.125
$ws.Range("M7").Value = -118.7
$ws.Range("N7").Value = -725.125
$ws.Range("H58").Value = 1884.3182
$ws.Range("I58").Value = 979.5
$ws.Range("J58").Value = 2970.1
$ws.Range("K58").Value = 979.5
$ws.Range("L58").Value = 2970.1
$ws.Range("M58").Value = -776.5
$ws.Range("N58").Value = -3376.1
$ws.Range("H94").Value = 1328.4
$ws.Range("I94").Value = 1111.6
$ws.Range("K94").Value = 1111.6
$ws.Range("M94").Value = -660.5999999999999
$ws.Range("H99").Value = 9392.706
$ws.Range("I99").Value = 10076.714
$ws.Range("K99").Value = 10076.714
$ws.Range("M99").Value = -8578.714
$ws.Range("H105").Value = 12622.454
$ws.Range("I105").Value = 1949.4
$ws.Range("K105").Value = 1949.4
$ws.Range("M105").Value = -202.4000000000001
$ws.Range("H126").Value = 9392.706
$ws.Range("I126").Value = 10076.714
$ws.Range("K126").Value = 30230.142
$ws.Range("M126").Value = -27760.142
$ws.Range("H132").Value = 48776.387
$ws.Range("I132").Value = 55523
$ws.Range("K132").Value = 166569
$ws.Range("M132").Value = -164039
$ws.Range("H136").Value = 1884.3182
$ws.Range("I136").Value = 979.5
$ws.Range("J136").Value = 2970.1
$ws.Range("K136").Value = 2938.5
$ws.Range("L136").Value = 8910.299999999999
$ws.Range("M136").Value = -388.5
$ws.Range("N136").Value = -14010.3
$ws.Range("H141").Value = 300771.1
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 300771.1
$ws.Range("K141").Value = 0
$ws.Range("M141").ClearContents()
$ws.Range("N141").Value = -311131.1

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H47").Value = 2000
$ws.Range("I47").Value = 1500
$ws.Range("J47").Value = 2250
$ws.Range("K47").Value = 4500
$ws.Range("L47").Value = 6750
$ws.Range("M47").Value = -4069
$ws.Range("N47").Value = -7612
$ws.Range("H51").Value = 1103.375
$ws.Range("I51").Value = 365.4
$ws.Range("J51").Value = 2333.3333
$ws.Range("K51").Value = 1096.2
$ws.Range("L51").Value = 6999.999899999999
$ws.Range("M51").Value = -636.1999999999998
$ws.Range("N51").Value = -7919.999899999999
$ws.Range("H56").Value = 16431.117
$ws.Range("I56").Value = 16431.117
$ws.Range("K56").Value = 16431.117
$ws.Range("M56").Value = -15901.117
$ws.Range("H59").Value = 1557.2858
$ws.Range("I59").Value = 1380.2
$ws.Range("J59").Value = 2000
$ws.Range("K59").Value = 4140.6
$ws.Range("L59").Value = 6000
$ws.Range("M59").Value = -3600.6
$ws.Range("N59").Value = -7080
$ws.Range("H80").Value = 992
$ws.Range("J80").Value = 992
$ws.Range("L80").Value = 2976
$ws.Range("N80").Value = -4848
$ws.Range("H83").Value = 992
$ws.Range("J83").Value = 992
$ws.Range("L83").Value = 8928
$ws.Range("N83").Value = -18288
$ws.Range("H113").Value = 3165.4119
$ws.Range("J113").Value = 3610.6667
$ws.Range("L113").Value = 10832.0001
$ws.Range("N113").Value = -15172.0001
$ws.Range("H132").Value = 4749.75
$ws.Range("I132").Value = 2999.5
$ws.Range("J132").Value = 6500
$ws.Range("K132").Value = 26995.5
$ws.Range("L132").Value = 58500
$ws.Range("M132").Value = -24465.5
$ws.Range("N132").Value = -63560

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 207950
$ws.Range("J70").Value = 3497
$ws.Range("L70").Value = 3497
$ws.Range("N70").Value = -4037
$ws.Range("H73").Value = 207950
$ws.Range("J73").Value = 3497
$ws.Range("L73").Value = 3497
$ws.Range("N73").Value = -5369
$ws.Range("H80").Value = 2200.3125
$ws.Range("I80").Value = 1867.409
$ws.Range("J80").Value = 2932.7
$ws.Range("K80").Value = 1867.409
$ws.Range("L80").Value = 2932.7
$ws.Range("M80").Value = -869.4090000000001
$ws.Range("N80").Value = -4928.7
$ws.Range("H83").Value = 2200.3125
$ws.Range("I83").Value = 1867.409
$ws.Range("J83").Value = 2932.7
$ws.Range("K83").Value = 9337.045
$ws.Range("L83").Value = 14663.5
$ws.Range("M83").Value = -4345.045
$ws.Range("N83").Value = -24647.5
$ws.Range("H102").Value = 3596.05
$ws.Range("I102").Value = 2059.4443
$ws.Range("K102").Value = 2059.4443
$ws.Range("M102").Value = -437.4443000000001
$ws.Range("H113").Value = 2831.5908
$ws.Range("I113").Value = 754.25
$ws.Range("K113").Value = 754.25
$ws.Range("M113").Value = 1415.75
$ws.Range("H122").Value = 1321.4736
$ws.Range("I122").Value = 1244.75
$ws.Range("K122").Value = 3734.25
$ws.Range("M122").Value = -1284.25
$ws.Range("H132").Value = 3851.9333
$ws.Range("I132").Value = 4064.9167
$ws.Range("K132").Value = 12194.7501
$ws.Range("M132").Value = -9664.750100000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1525.3948
$ws.Range("I46").Value = 639.5357
$ws.Range("J46").Value = 4005.8
$ws.Range("K46").Value = 639.5357
$ws.Range("L46").Value = 4005.8
$ws.Range("M46").Value = -451.5357
$ws.Range("N46").Value = -4381.8
$ws.Range("H55").Value = 691.5333000000001
$ws.Range("I55").Value = 237.16667
$ws.Range("J55").Value = 994.44446
$ws.Range("K55").Value = 237.16667
$ws.Range("L55").Value = 994.44446
$ws.Range("M55").Value = -64.16667000000001
$ws.Range("N55").Value = -1340.44446
$ws.Range("H63").Value = 68332.664
$ws.Range("J63").Value = 65998
$ws.Range("L63").Value = 65998
$ws.Range("N63").Value = -67496
$ws.Range("H66").Value = 68332.664
$ws.Range("J66").Value = 65998
$ws.Range("L66").Value = 197994
$ws.Range("N66").Value = -205482
$ws.Range("H93").Value = 2047.875
$ws.Range("I93").Value = 1122.8182
$ws.Range("K93").Value = 1122.8182
$ws.Range("M93").Value = 125.1818000000001
$ws.Range("H100").Value = 3011.0908
$ws.Range("I100").Value = 2686.158
$ws.Range("K100").Value = 2686.158
$ws.Range("M100").Value = -2145.158
$ws.Range("H122").Value = 4141.9707
$ws.Range("I122").Value = 2486.923
$ws.Range("J122").Value = 5166.524
$ws.Range("K122").Value = 7460.768999999999
$ws.Range("L122").Value = 15499.572
$ws.Range("M122").Value = -5010.768999999999
$ws.Range("N122").Value = -20399.572
$ws.Range("H132").Value = 1702.9584
$ws.Range("I132").Value = 1743.1
$ws.Range("K132").Value = 5229.299999999999
$ws.Range("M132").Value = -2699.299999999999
$ws.Range("H136").Value = 2432.3
$ws.Range("I136").Value = 1970.5834
$ws.Range("K136").Value = 5911.7502
$ws.Range("M136").Value = -3361.7502

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 76924380
$ws.Range("I100").Value = 90910420
$ws.Range("J100").Value = 1250
$ws.Range("K100").Value = 181820840
$ws.Range("L100").Value = 2500
$ws.Range("M100").Value = -181820299
$ws.Range("N100").Value = -3582
$ws.Range("H107").Value = 499.25
$ws.Range("I107").Value = 499.14285
$ws.Range("K107").Value = 1497.42855
$ws.Range("M107").Value = 422.5714499999999
$ws.Range("H113").Value = 684.1
$ws.Range("I113").Value = 395.5625
$ws.Range("J113").Value = 1013.8571
$ws.Range("K113").Value = 1186.6875
$ws.Range("L113").Value = 3041.5713
$ws.Range("M113").Value = 983.3125
$ws.Range("N113").Value = -7381.5713
$ws.Range("H122").Value = 1465.1875
$ws.Range("I122").Value = 1224.7037
$ws.Range("K122").Value = 3674.1111
$ws.Range("M122").Value = -1224.1111
$ws.Range("H136").Value = 1192.3256
$ws.Range("I136").Value = 796.17645
$ws.Range("K136").Value = 2388.52935
$ws.Range("M136").Value = 161.4706499999998
